# Apply updated crypto price/volume data per Mon May 20 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.881.02"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.085.48"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.76%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.085.02"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "3.598.78"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "66.789.42"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "3.082.92"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "484.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.684"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "0.0₃0927"
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.941"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("B37").Value = "Arweave"
$ws.Range("C37").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "46.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.311"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "48.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.71%  "
$ws.Range("D45").Value = "2.778.96"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "366.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("E51").Value = "  +5.30%  "
